$d = $word.ActiveDocument

$d.Content.Find.Execute("502÷5=100, 2", $true, $false, $false, $false, $false, $true, 1, $false, "221÷4=55, 1", 2)
$d.Content.Find.Execute("669÷5=133, 4", $true, $false, $false, $false, $false, $true, 1, $false, "122÷2=61, 0", 2)
$d.Content.Find.Execute("420÷4=105, 0", $true, $false, $false, $false, $false, $true, 1, $false, "586÷9=65, 1", 2)
$d.Content.Find.Execute("875÷4=218, 3", $true, $false, $false, $false, $false, $true, 1, $false, "151÷4=37, 3", 2)
$d.Content.Find.Execute("590÷6=98, 2", $true, $false, $false, $false, $false, $true, 1, $false, "481÷6=80, 1", 2)
$d.Content.Find.Execute("585÷7=83, 4", $true, $false, $false, $false, $false, $true, 1, $false, "535÷2=267, 1", 2)
$d.Content.Find.Execute("554÷4=138, 2", $true, $false, $false, $false, $false, $true, 1, $false, "324÷2=162, 0", 2)
$d.Content.Find.Execute("381÷4=95, 1", $true, $false, $false, $false, $false, $true, 1, $false, "410÷4=102, 2", 2)
$d.Content.Find.Execute("979÷9=108, 7", $true, $false, $false, $false, $false, $true, 1, $false, "245÷6=40, 5", 2)
$d.Content.Find.Execute("106÷5=21, 1", $true, $false, $false, $false, $false, $true, 1, $false, "843÷8=105, 3", 2)
$d.Content.Find.Execute("631÷3=210, 1", $true, $false, $false, $false, $false, $true, 1, $false, "702÷4=175, 2", 2)
$d.Content.Find.Execute("982÷7=140, 2", $true, $false, $false, $false, $false, $true, 1, $false, "103÷2=51, 1", 2)
$d.Content.Find.Execute("937÷4=234, 1", $true, $false, $false, $false, $false, $true, 1, $false, "853÷2=426, 1", 2)
$d.Content.Find.Execute("645÷2=322, 1", $true, $false, $false, $false, $false, $true, 1, $false, "529÷7=75, 4", 2)
$d.Content.Find.Execute("888÷8=111, 0", $true, $false, $false, $false, $false, $true, 1, $false, "127÷5=25, 2", 2)
$d.Content.Find.Execute("274÷4=68, 2", $true, $false, $false, $false, $false, $true, 1, $false, "718÷9=79, 7", 2)
$d.Content.Find.Execute("370÷9=41, 1", $true, $false, $false, $false, $false, $true, 1, $false, "255÷8=31, 7", 2)
$d.Content.Find.Execute("923÷8=115, 3", $true, $false, $false, $false, $false, $true, 1, $false, "165÷6=27, 3", 2)
$d.Content.Find.Execute("773÷8=96, 5", $true, $false, $false, $false, $false, $true, 1, $false, "812÷7=116, 0", 2)
$d.Content.Find.Execute("991÷5=198, 1", $true, $false, $false, $false, $false, $true, 1, $false, "468÷2=234, 0", 2)
$d.Content.Find.Execute("533÷7=76, 1", $true, $false, $false, $false, $false, $true, 1, $false, "488÷3=162, 2", 2)
$d.Content.Find.Execute("536÷7=76, 4", $true, $false, $false, $false, $false, $true, 1, $false, "278÷9=30, 8", 2)
$d.Content.Find.Execute("525÷5=105, 0", $true, $false, $false, $false, $false, $true, 1, $false, "128÷9=14, 2", 2)
$d.Content.Find.Execute("949÷6=158, 1", $true, $false, $false, $false, $false, $true, 1, $false, "819÷4=204, 3", 2)
